# Lvl Meter switched to moving average - Currently no Automatic Gain Control
#
# The source sheet (Sheet1) is a small calculator: B16/B18/F22 are the only
# true "inputs" in the affected region; every other changed cell (B19, B20,
# G22, G23, B24, G24, B26, B27) is a formula that recalculates automatically
# once those inputs change. We just have to set the new input values and
# move the active selection to match what the author left selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Input cells that actually changed -------------------------------------
# "Fs ref" (B16): 8.192 -> 4.096
$ws.Range("B16").Value = 4.096

# "N" (B18, the 0x25 div reference count): 2048 -> 4096
$ws.Range("B18").Value = 4096

# "0x25 div" (F22): 16 -> 8
$ws.Range("F22").Value = 8

# --- Selection / viewport ----------------------------------------------------
# Author's cursor ended up on B17 (was B26), with the view scrolled so row 4
# is the first visible row.
$ws.Range("B17").Select()

# Best-effort: try to scroll the window so A4 is the top-left visible cell,
# matching sheetView's topLeftCell="A4". Harmless if unsupported.
try { $excel.ActiveWindow.ScrollRow = 4 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
